# Apply updates to the "展览" and "全部类型" sheets.
# Both sheets carry identical data tables, and both receive the same edits.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Row 2: F unchanged; G2 0 -> "不可售"
    $ws.Range("G2").Value = "不可售"

    # Row 3: F3 7740 -> 7766
    $ws.Range("F3").Value = 7766

    # Row 4: F4 530 -> 531; G4 0 -> "已售罄"
    $ws.Range("F4").Value = 531
    $ws.Range("G4").Value = "已售罄"

    # Row 5: F5 326 -> 327; G5 258 -> "已售罄"
    $ws.Range("F5").Value = 327
    $ws.Range("G5").Value = "已售罄"

    # Row 6: F6 41 -> 43; G6 55 -> 65
    $ws.Range("F6").Value = 43
    $ws.Range("G6").Value = 65

    # Row 7: F7 22 -> 23
    $ws.Range("F7").Value = 23

    # Row 9: F9 5986 -> 6096
    $ws.Range("F9").Value = 6096

    # Row 10: F10 153 -> 162
    $ws.Range("F10").Value = 162

    # Row 13: F13 1820 -> 1840
    $ws.Range("F13").Value = 1840

    # Row 14: F14 1337 -> 1376
    $ws.Range("F14").Value = 1376

    # Row 15: F15 285 -> 289
    $ws.Range("F15").Value = 289

    # Row 16: F16 629 -> 770
    $ws.Range("F16").Value = 770

    # Row 17: F17 147 -> 174; G17 0 -> "不可售"
    $ws.Range("F17").Value = 174
    $ws.Range("G17").Value = "不可售"

    # Row 18: F18 5539 -> 5547
    $ws.Range("F18").Value = 5547

    # Row 19: F19 66 -> 67
    $ws.Range("F19").Value = 67
}
